$d = $word.ActiveDocument

# Locate the insertion point: right after "...on his first day" and
# before " he took it upon himself..." — a comma needs to be inserted there.
$r = $d.Content
$found = $r.Find.Execute("first day", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor text 'first day'"
}

# Collapse the found range to its end (right after "day") and type the comma.
$r.Collapse(0)
$r.InsertAfter(",")

# The comma now occupies the single character just inserted; remember where
# it ends so we can drop the _GoBack bookmark there, mirroring Word's
# behaviour of moving _GoBack to the site of the most recent edit.
$commaEnd = $r.Start + 1
$gobackRange = $d.Range($commaEnd, $commaEnd)

# Remove the _GoBack bookmark from its old location (end of document) and
# recreate it at the new edit location.
$old = $d.Bookmarks("_GoBack")
$old.Delete()
$d.Bookmarks.Add("_GoBack", $gobackRange)
